$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new data point for 2026/01/22 (Thu) was recorded, so a row is inserted
# right after the existing block of 2026/01/22 rows (before the 2026/12/29
# block), shifting every row from 677 downward by one.
$ws.Rows.Item(677).Insert()

# Copy the row above (still the same 2026/01/22 / 木 day) into the freshly
# inserted blank row so that the date/weekday cells keep their original
# "text" storage (plain inline string, no number formatting) instead of
# Excel auto-converting a typed "2026/01/22" string into a date serial.
$ws.Range("A676:D676").Copy()
$ws.Range("A677:D677").PasteSpecial()

# Now overwrite the time/ranking columns with the new observation's values.
$ws.Cells.Item(677, 3).Value = 16
$ws.Cells.Item(677, 4).Value = 20
